$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.073.33"
$ws.Range("E2").Value = "  -2.22%  "

# Row 3
$ws.Range("D3").Value = "2.169.69"
$ws.Range("E3").Value = "  -2.13%  "

# Row 4
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "236.73"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.09%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.614"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.98%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "70.16"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -4.65%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  -5.19%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "40.37"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -6.87%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0925"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.40%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "54.47"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -5.19%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.101"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.34%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.76"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -5.11%  "

# Row 15
$ws.Range("D15").Value = "2.491.62"
$ws.Range("E15").Value = "  -2.24%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "13.95"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -2.32%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.180.59"
$ws.Range("E17").Value = "  -1.49%  "

# Row 18
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.802"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -4.82%  "

# Row 19
$ws.Range("D19").Value = "40.953.71"
$ws.Range("E19").Value = "  -2.22%  "

# Row 20
$ws.Range("E20").Value = "  -6.94%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "70.40"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.14%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.91"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.97%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.81"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -5.92%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "225.67"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.75%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.95"
$c.Style = "Normal"

# Row 26
$ws.Range("E26").Value = "  -0.09%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.86"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -5.68%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "3.53"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.38%  "

# Row 29
$ws.Range("E29").Value = "  -3.11%  "

# Row 30
$ws.Range("E30").Value = "  +0.90%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "167.58"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.22%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "19.91"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.34%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "30.82"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +5.55%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0764"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.91%  "

# Row 35
$ws.Range("E35").Value = "  -7.63%  "

# Row 36
$ws.Range("E36").Value = "  -3.31%  "

# Row 37
$ws.Range("E37").Value = "  -7.55%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.13"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.58%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0284"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -5.90%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "12.00"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -6.48%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.07"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.29%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.42"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.81%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "59.45"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -10.09%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.190"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -4.69%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0974"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.23%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "8.26"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -5.04%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "97.73"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -6.20%  "

# Row 48
$ws.Range("E48").Value = "  -2.51%  "

# Row 49
$ws.Range("E49").Value = "  -2.98%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.22"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -8.33%  "

# Row 51
$ws.Range("E51").Value = "  -2.85%  "
